$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.729.83'
$ws.Range('E2').Value = '  +3.13%  '
$ws.Range('D3').Value = '1.864.27'
$ws.Range('E3').Value = '  +2.93%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.037'
$ws.Range('E4').Value = '  +2.98%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '324.52'
$ws.Range('E5').Value = '  +3.86%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.033'
$ws.Range('E6').Value = '  +2.73%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4414'
$ws.Range('E7').Value = '  +2.85%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3799'
$ws.Range('E8').Value = '  +2.66%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07462'
$ws.Range('E9').Value = '  +3.12%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8843'
$ws.Range('E10').Value = '  +2.18%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '21.81'
$ws.Range('E11').Value = '  +2.29%  '
$ws.Range('D12').Value = '1.890.81'
$ws.Range('E12').Value = '  -7.57%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.554'
$ws.Range('E13').Value = '  +2.87%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.751'
$ws.Range('E14').Value = '  +1.80%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.07219'
$ws.Range('E15').Value = '  +4.09%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '83.89'
$ws.Range('E16').Value = '  +3.81%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.038'
$ws.Range('E17').Value = '  +2.60%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000009104'
$ws.Range('E18').Value = '  +2.16%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.033'
$ws.Range('E19').Value = '  +2.74%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '15.44'
$ws.Range('E20').Value = '  +1.65%  '
$ws.Range('D21').Value = '27.757.72'
$ws.Range('E21').Value = '  +3.10%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.310'
$ws.Range('E22').Value = '  +2.21%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.43'
$ws.Range('E23').Value = '  +4.36%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.955'
$ws.Range('E24').Value = '  +3.72%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '158.10'
$ws.Range('E25').Value = '  +2.36%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '18.87'
$ws.Range('E26').Value = '  +2.98%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.001'
$ws.Range('E27').Value = '  +4.14%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '5.306'
$ws.Range('E28').Value = '  +1.22%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '117.53'
$ws.Range('E29').Value = '  +2.54%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.09089'
$ws.Range('E30').Value = '  +1.51%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.214'
$ws.Range('E31').Value = '  +4.89%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.7695'
$ws.Range('E32').Value = '  +3.47%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.064'
$ws.Range('E33').Value = '  +9.24%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.579'
$ws.Range('E34').Value = '  +3.30%  '
$ws.Range('E35').Value = '  +2.77%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.164'
$ws.Range('E36').Value = '  +3.59%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.01987'
$ws.Range('E37').Value = '  +3.15%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.05349'
$ws.Range('E38').Value = '  +2.23%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.5195'
$ws.Range('E39').Value = '  +1.97%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.838'
$ws.Range('E40').Value = '  +3.28%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.1693'
$ws.Range('E41').Value = '  +2.48%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.846'
$ws.Range('E42').Value = '  +5.56%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.706'
$ws.Range('E43').Value = '  +5.01%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '109.56'
$ws.Range('E44').Value = '  +1.96%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '10.62'
$ws.Range('E45').Value = '  +2.24%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.725'
$ws.Range('E46').Value = '  +4.65%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4696'
$ws.Range('E47').Value = '  +3.03%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.06428'
$ws.Range('E48').Value = '  +2.59%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.872'
$ws.Range('E49').Value = '  +3.99%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '39.86'
$ws.Range('E50').Value = '  +4.51%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '64.49'
$ws.Range('E51').Value = '  +2.03%  '
